# Generate Report for Handback
# Row 3 (f1607909-...md) moves from "Ready for handoff" to
# "Handed back: in sync with en-US" on all three sheets, with fresh
# handback timestamps and the stale "not latest" error cleared out.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet: row 3 status columns (zh-cn / de-de) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText

# --- zh-cn sheet: row 3 ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $statusText
$zhcn.Range("K3").Value = "2016-08-19 04:43:30"
$zhcn.Range("P3").Value = ""

# --- de-de sheet: row 3 ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $statusText
$dede.Range("K3").Value = "2016-08-19 04:43:37"
$dede.Range("P3").Value = ""
